$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 175, pushing the existing rows 175-195 down to 176-196
$ws.Rows("175:175").Insert()

# Populate the new weekly data row (Santina / Primera, week of 2023-12-05)
$ws.Range("A175").Value = 7
$ws.Range("B175").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C175").Value = "Ñuble"
$ws.Range("D175").Value = 45265
$ws.Range("E175").Value = 16
$ws.Range("F175").Value = "Fruta"
$ws.Range("G175").Value = 100103
$ws.Range("H175").Value = "Frutos de hueso (carozo)"
$ws.Range("I175").Value = 100103001
$ws.Range("J175").Value = "Cereza"
$ws.Range("K175").Value = "Santina"
$ws.Range("L175").Value = "Primera"
$ws.Range("M175").Value = 150
$ws.Range("N175").Value = 10000
$ws.Range("O175").Value = 10000
$ws.Range("P175").Value = 10000
$ws.Range("Q175").Value = "$/bandeja 10 kilos"
$ws.Range("R175").Value = "Provincia de Curicó"
$ws.Range("S175").Value = 1000
$ws.Range("T175").Value = 10
